$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 827 (shifts the existing row 827..868 down to 828..869)
$ws.Rows.Item(827).Insert()

# Fill the newly inserted row with its data. Column A holds dates stored as
# plain text (matching the rest of the sheet), so force text formatting
# before assigning the value and then clear the format so no extra style
# index is left behind on the cell.
$ws.Cells.Item(827, 1).NumberFormat = "@"
$ws.Cells.Item(827, 1).Value = "2026/02/21"
$ws.Cells.Item(827, 1).ClearFormats()

$ws.Cells.Item(827, 2).Value = "土"
$ws.Cells.Item(827, 3).Value = 16
$ws.Cells.Item(827, 4).Value = 32
